$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44224
$ws.Range("L2").Value = "Primera"
$ws.Range("N2").Value = 16500
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16750
$ws.Range("S2").Value = 931

# Row 3
$ws.Range("D3").Value = 44224
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 13500
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 13750
$ws.Range("S3").Value = 764

# Row 4
$ws.Range("D4").Value = 44335
$ws.Range("L4").Value = "Especial"
$ws.Range("N4").Value = 19500
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19750
$ws.Range("S4").Value = 1097

# Row 5
$ws.Range("D5").Value = 44335
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 17500
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17750
$ws.Range("S5").Value = 986

# Row 6
$ws.Range("D6").Value = 44335
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 12500
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12750
$ws.Range("S6").Value = 708

# Row 7
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 12500
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12750
$ws.Range("S7").Value = 708

# Row 8
$ws.Range("D8").Value = 44272
$ws.Range("L8").Value = "Primera"
$ws.Range("N8").Value = 10500
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 10750
$ws.Range("S8").Value = 597

# Row 9
$ws.Range("D9").Value = 44272
$ws.Range("L9").Value = "Segunda"
$ws.Range("N9").Value = 8500
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 8750
$ws.Range("S9").Value = 486

# Row 10
$ws.Range("D10").Value = 44294
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 360
$ws.Range("N10").Value = 12500
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 12750
$ws.Range("S10").Value = 708

# Row 11
$ws.Range("D11").Value = 44294
$ws.Range("M11").Value = 240
$ws.Range("N11").Value = 10500
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 10750
$ws.Range("S11").Value = 597

# Row 12
$ws.Range("D12").Value = 44294
$ws.Range("M12").Value = 240
$ws.Range("N12").Value = 8500
$ws.Range("O12").Value = 9000
$ws.Range("P12").Value = 8750
$ws.Range("S12").Value = 486

# Row 13
$ws.Range("D13").Value = 44385
$ws.Range("M13").Value = 120
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14500
$ws.Range("P13").Value = 14250
$ws.Range("S13").Value = 792

# Row 14
$ws.Range("D14").Value = 44385
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 11000
$ws.Range("O14").Value = 11500
$ws.Range("P14").Value = 11250
$ws.Range("S14").Value = 625

# Row 15
$ws.Range("D15").Value = 44385
$ws.Range("M15").Value = 240
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8500
$ws.Range("P15").Value = 8250
$ws.Range("S15").Value = 458

# Row 16
$ws.Range("D16").Value = 44385
$ws.Range("L16").Value = "Tercera"
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 5000
$ws.Range("O16").Value = 5500
$ws.Range("P16").Value = 5250
$ws.Range("S16").Value = 292

# Row 17
$ws.Range("D17").Value = 44308
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 15500
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15750
$ws.Range("S17").Value = 875

# Row 18
$ws.Range("D18").Value = 44308
$ws.Range("L18").Value = "Primera"
$ws.Range("N18").Value = 13500
$ws.Range("O18").Value = 14000
$ws.Range("P18").Value = 13750
$ws.Range("S18").Value = 764

# Row 19
$ws.Range("D19").Value = 44308
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 10500
$ws.Range("O19").Value = 11000
$ws.Range("P19").Value = 10750
$ws.Range("S19").Value = 597

# Row 20
$ws.Range("D20").Value = 44279
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 12500
$ws.Range("O20").Value = 13000
$ws.Range("P20").Value = 12750
$ws.Range("S20").Value = 708

# Row 21
$ws.Range("D21").Value = 44279
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 240
$ws.Range("N21").Value = 10500
$ws.Range("O21").Value = 11000
$ws.Range("P21").Value = 10750
$ws.Range("S21").Value = 597

# Row 22
$ws.Range("D22").Value = 44279
$ws.Range("L22").Value = "Segunda"
$ws.Range("N22").Value = 8500
$ws.Range("O22").Value = 9000
$ws.Range("P22").Value = 8750
$ws.Range("S22").Value = 486

# Row 23
$ws.Range("L23").Value = "Especial"
$ws.Range("N23").Value = 12500
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 12750
$ws.Range("S23").Value = 708

# Row 24
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 240
$ws.Range("N24").Value = 10500
$ws.Range("O24").Value = 11000
$ws.Range("P24").Value = 10750
$ws.Range("S24").Value = 597

# Row 25
$ws.Range("D25").Value = 44280
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 8500
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 8750
$ws.Range("S25").Value = 486

# Row 26
$ws.Range("D26").Value = 44293
$ws.Range("L26").Value = "Especial"
$ws.Range("M26").Value = 400
$ws.Range("N26").Value = 12500
$ws.Range("O26").Value = 13000
$ws.Range("P26").Value = 12750
$ws.Range("S26").Value = 708

# Row 27
$ws.Range("D27").Value = 44293
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 508
$ws.Range("N27").Value = 10500
$ws.Range("O27").Value = 11000
$ws.Range("P27").Value = 10746
$ws.Range("S27").Value = 597

# Row 28
$ws.Range("D28").Value = 44293
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 400
$ws.Range("N28").Value = 8500
$ws.Range("O28").Value = 9000
$ws.Range("P28").Value = 8750
$ws.Range("S28").Value = 486

# Row 29
$ws.Range("D29").Value = 44384
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 160
$ws.Range("N29").Value = 14500
$ws.Range("O29").Value = 15000
$ws.Range("P29").Value = 14750
$ws.Range("S29").Value = 819

# Row 30
$ws.Range("D30").Value = 44384
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 500
$ws.Range("N30").Value = 11500
$ws.Range("O30").Value = 12000
$ws.Range("P30").Value = 11750
$ws.Range("S30").Value = 653

# Row 31
$ws.Range("D31").Value = 44384
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 400
$ws.Range("N31").Value = 8500
$ws.Range("O31").Value = 9000
$ws.Range("P31").Value = 8750
$ws.Range("S31").Value = 486

# Row 32
$ws.Range("L32").Value = "Tercera"
$ws.Range("M32").Value = 300
$ws.Range("N32").Value = 5500
$ws.Range("O32").Value = 6000
$ws.Range("P32").Value = 5750
$ws.Range("S32").Value = 319

# Row 33
$ws.Range("D33").Value = 44286
$ws.Range("L33").Value = "Especial"
$ws.Range("M33").Value = 700
$ws.Range("N33").Value = 12500
$ws.Range("O33").Value = 13000
$ws.Range("P33").Value = 12750
$ws.Range("S33").Value = 708

# Row 34
$ws.Range("D34").Value = 44286
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 500
$ws.Range("N34").Value = 10500
$ws.Range("O34").Value = 11000
$ws.Range("P34").Value = 10750
$ws.Range("S34").Value = 597

# Row 35
$ws.Range("D35").Value = 44286
$ws.Range("L35").Value = "Segunda"
$ws.Range("N35").Value = 8500
$ws.Range("O35").Value = 9000
$ws.Range("P35").Value = 8750
$ws.Range("S35").Value = 486

# Row 36
$ws.Range("D36").Value = 44273
$ws.Range("M36").Value = 160

# Row 37
$ws.Range("D37").Value = 44273
$ws.Range("M37").Value = 240

# Row 38
$ws.Range("D38").Value = 44273
$ws.Range("M38").Value = 200
